$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

$ws.Range("A31").Value = "Rdata "
$ws.Range("B31").Value = "2022-06-14 16-10-45_chunks"
$ws.Range("C31").Value = "chunking Kriegstein"
$ws.Range("D31").Value = "SCTv2 corrected pipeline rechunking Kriegstein ref data"
$ws.Range("F31").Value = "rerun SCTv2 corrected pipeline"

$ws.Range("F31").Select()
